# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Simple numeric (Client/Total Samples count) updates
$ws.Range("B15").Value = 63373
$ws.Range("B16").Value = 69345
$ws.Range("B21").Value = 276524
$ws.Range("B22").Value = 637128

# Row 25/26 swap adapter-driver text + counts
$ws.Range("A25").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1"
$ws.Range("B25").Value = 69578
$ws.Range("A26").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.250.10.1"
$ws.Range("B26").Value = 78477

$ws.Range("B27").Value = 338880
$ws.Range("B29").Value = 459268
$ws.Range("B30").Value = 143869
$ws.Range("B32").Value = 31330

# Row 33/34 swap adapter-driver text + counts
$ws.Range("A33").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4"
$ws.Range("B33").Value = 96526
$ws.Range("A34").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.160.0.4"
$ws.Range("B34").Value = 3654692

$ws.Range("B35").Value = 158283
$ws.Range("B37").Value = 106139
$ws.Range("B39").Value = 109293
$ws.Range("B40").Value = 101951

$ws.Range("B43").Value = 191994
$ws.Range("D43").Value = 99.90000000000001

$ws.Range("B44").Value = 244856
$ws.Range("B46").Value = 325504
$ws.Range("B47").Value = 68450
$ws.Range("B49").Value = 96727
$ws.Range("B50").Value = 122671
$ws.Range("B52").Value = 689912
$ws.Range("B54").Value = 211798
$ws.Range("B55").Value = 69430
$ws.Range("B57").Value = 310711
$ws.Range("B58").Value = 90508

# Row 59/60 swap adapter-driver text + counts
$ws.Range("A59").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.0.0.5"
$ws.Range("B59").Value = 196098
$ws.Range("A60").Value = "Intel(R) Wireless-AC 9560 160MHz - 21.40.2.2"
$ws.Range("B60").Value = 140512

$ws.Range("B64").Value = 451638
$ws.Range("B65").Value = 52515
$ws.Range("B66").Value = 116738
$ws.Range("B68").Value = 73817
